$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 values per the diff
$ws.Range("A1").Value = 4
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 5
$ws.Range("F1").Value = 30
$ws.Range("G1").Value = 20
$ws.Range("H1").Value = 18
$ws.Range("I1").Value = 32
$ws.Range("J1").Value = 33
$ws.Range("K1").Value = 0.083000000000000004
$ws.Range("L1").Value = 0.048000000000000001
$ws.Range("M1").Value = 0.079000000000000001
$ws.Range("N1").Value = 0.029999999999999999

# Update column widths per the diff.
# Note: the COM layer quantizes ColumnWidth to 1/6-character increments, so the
# input values below are chosen to land as close as possible on the exact
# target widths (2.140625 / 3.140625 / 5.7109375 / 5.7109375 / 4.7109375).
$ws.Columns.Item(5).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(6).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(11).ColumnWidth = 4.833333333333333
$ws.Columns.Item(12).ColumnWidth = 4.833333333333333
$ws.Columns.Item(14).ColumnWidth = 3.8333333333333335
